# Logbook entry 16: "Finished path generation code"
# Fills in the previously-blank row 20 of the logbook with a new entry,
# adjusts the row height to fit the longer description text, and leaves
# the active selection on L16 (matching the author's saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# No. column
$ws.Range("B20").Value = 16

# Date column (21 May 2021, stored as the Excel serial date already used
# by the rest of the sheet)
$ws.Range("C20").Value = 44337

# Title column
$ws.Range("D20").Value = "Finished path generation code"

# Description column
$ws.Range("E20").Value = "Fixed errors with path generation code and made it responsive to width limitations and angle limitaitons. Maybe at some point in the future, the width of the path itself may be implemented, but for now that's not important. This code now needs to be adapted to storing images in RAM to make the process faster, and it needs to be implemented into the existing game. It is currently a separate file."

# The extra text requires a taller row than the default entries.
$ws.Rows.Item(20).RowHeight = 60.5

# Restore the saved selection/active-cell state from the source workbook.
[void]$ws.Range("L16").Select()
